$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4298446476459503
$ws.Range("B1").Value = 0.8415980339050293
$ws.Range("C1").Value = 4.973837375640869
$ws.Range("D1").Value = 2.160943269729614
$ws.Range("E1").Value = 1.301877737045288
